$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceOne = 1, wdReplaceAll = 2, wdReplaceNone = 0
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,
#              MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)

# --- Change 1: opening paragraph - rewrite the two closing sentences ---
$old1 = "make an impact during their Rookie season (i.e. their first season) in the NFL. It is our goal to provide a tool with which NFL teams can rank potential running back choices and pick one that will be able to have the most impact during their first year."
$new1 = "be considered worth being drafted. It is our goal to provide a tool with which NFL teams can determine which players deserve extra attention by their scouting offices."
$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2: "Rookie running back data" -> "NFL running back data" ---
$old2 = "the ones with Rookie running back data from 2010"
$new2 = "the ones with NFL running back data from 2010"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Change 3: swap the NFL.com stats hyperlink display text/URL ---
$old3 = "http://www.nfl.com/stats/categorystats?seasonType=REG&d-447263-n=1&d-447263-o=2&d-447263-p=1&statisticPositionCategory=RUNNING_BACK&d-447263-s=RUSHING_YARDS&tabSeq=1&season=2016&Submit=Go&experience=0&archive=false&conference=null&qualified=true"
$new3 = "http://www.nfl.com/stats/categorystats?archive=true&conference=null&statisticCategory=RUSHING&season=2016&seasonType=REG&experience=&tabSeq=0&qualified=false&Submit=Go"
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# --- Change 4a: "NFL Rookie information" -> "NFL player information" ---
$old4a = "our NFL Rookie information from the corresponding years"
$new4a = "our NFL player information from the corresponding years"
$d.Content.Find.Execute($old4a, $true, $false, $false, $false, $false, $true, 1, $false, $new4a, 2) | Out-Null

# --- Change 4b: "allow us to rank a set of" -> "allow us to pass judgment on a set of" ---
$old4b = "create a model which will allow us to rank a set of"
$new4b = "create a model which will allow us to pass judgment on a set of"
$d.Content.Find.Execute($old4b, $true, $false, $false, $false, $false, $true, 1, $false, $new4b, 2) | Out-Null

# --- Change 5: closing paragraph rewording ---
$old5 = "rank their ability to have an impact during their Rookie season."
$new5 = "determine if they are worth scouting for the upcoming NFL draft."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

# append the extra trailing space that now appears at the very end of the last paragraph
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRng = $lastPara.Range
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1) | Out-Null
$endRng.Collapse(0)
$endRng.InsertAfter(" ")

# --- Move the _GoBack bookmark from the end of the last paragraph to the end of the
#     "incoming college running backs." paragraph (the one just before it) ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$findRng = $d.Content
$findRng.Find.Execute("incoming college running backs. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $findRng)

Write-Output "done"
